$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 10: Inscritos 369 -> 370
$ws.Range("E10").Value = 370

# Row 21: Inscritos 116 -> 117, Pagos 59 -> 60, Inscricoes homologadas 59 -> 60
$ws.Range("E21").Value = 117
$ws.Range("F21").Value = 60
$ws.Range("H21").Value = 60

# Row 22: Inscritos 136 -> 137
$ws.Range("E22").Value = 137

# Row 40: Inscritos 207 -> 208
$ws.Range("E40").Value = 208

# Row 44: Inscritos 246 -> 247
$ws.Range("E44").Value = 247

# Row 47: Inscritos 343 -> 344, Pagos 158 -> 159, Inscricoes homologadas 158 -> 159
$ws.Range("E47").Value = 344
$ws.Range("F47").Value = 159
$ws.Range("H47").Value = 159

# Row 49: Inscritos 223 -> 224, Pagos 93 -> 94, Inscricoes homologadas 93 -> 94
$ws.Range("E49").Value = 224
$ws.Range("F49").Value = 94
$ws.Range("H49").Value = 94
